$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for c18ffc4e... row (Overview row 3)
# also shared with "Correspond Handoff Datetime" on de-de row 3 (same original value)
$wsOverview.Range("G3").Value = "2016-08-30 08:28:40"
$wsDeDe.Range("H3").Value = "2016-08-30 08:28:40"

# zh-cn row 3 (c18ffc4e...): Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H3").Value = "2016-08-30 08:28:27"
$wsZhCn.Range("K3").Value = "2016-08-30 08:29:26"

# de-de row 3 (c18ffc4e...): Correspond Handback DateTime
$wsDeDe.Range("K3").Value = "2016-08-30 08:29:46"
